# Update "paises.xlsx" - countries & provincias Spain refresh
# Applies the data refresh captured in the commit "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: "Datos actualizados a 15 de Octubre de 2020 a las 11:35" -> "...12:52" ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 12:52"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 8153740
$ws.Cells.Item(4, 3).Value = 3697
$ws.Cells.Item(4, 5).Value = 2652217
$ws.Cells.Item(4, 7).Value = 29
$ws.Cells.Item(4, 8).Value = 221872

# --- Row 16: Iran ---
$ws.Cells.Item(16, 2).Value = 517835
$ws.Cells.Item(16, 3).Value = 4616
$ws.Cells.Item(16, 4).Value = 418054
$ws.Cells.Item(16, 5).Value = 70176
$ws.Cells.Item(16, 7).Value = 256
$ws.Cells.Item(16, 8).Value = 29605

# --- Row 32: Rumania ---
$ws.Cells.Item(32, 2).Value = 168490
$ws.Cells.Item(32, 3).Value = 4013
$ws.Cells.Item(32, 4).Value = 125009
$ws.Cells.Item(32, 5).Value = 37807
$ws.Cells.Item(32, 7).Value = 73
$ws.Cells.Item(32, 8).Value = 5674

# --- Rows 40/41: Nepal and Republica Dominicana swap rank, Nepal gets fresh data ---
$ws.Cells.Item(40, 1).Value = "Nepal"
$ws.Cells.Item(40, 2).Value = 121745
$ws.Cells.Item(40, 3).Value = 3749
$ws.Cells.Item(40, 4).Value = 84518
$ws.Cells.Item(40, 5).Value = 36533
$ws.Cells.Item(40, 7).Value = 19
$ws.Cells.Item(40, 8).Value = 694

$ws.Cells.Item(41, 1).Value = "Republica Dominicana"
$ws.Cells.Item(41, 2).Value = 119662
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 95460
$ws.Cells.Item(41, 5).Value = 22016
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 2186

# --- Row 45: Oman ---
$ws.Cells.Item(45, 2).Value = 108296
$ws.Cells.Item(45, 3).Value = 520
$ws.Cells.Item(45, 4).Value = 94229
$ws.Cells.Item(45, 5).Value = 12996
$ws.Cells.Item(45, 7).Value = 10
$ws.Cells.Item(45, 8).Value = 1071

# --- Row 58: Suiza ---
$ws.Cells.Item(58, 2).Value = 71317
$ws.Cells.Item(58, 3).Value = 2613
$ws.Cells.Item(58, 5).Value = 19407
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 2110

# --- Row 82: El Salvador ---
$ws.Cells.Item(82, 2).Value = 31061
$ws.Cells.Item(82, 3).Value = 295
$ws.Cells.Item(82, 4).Value = 26311
$ws.Cells.Item(82, 5).Value = 3842
$ws.Cells.Item(82, 7).Value = 4
$ws.Cells.Item(82, 8).Value = 908

# --- Row 93: Malasia ---
$ws.Cells.Item(93, 2).Value = 18129
$ws.Cells.Item(93, 3).Value = 589
$ws.Cells.Item(93, 4).Value = 12014
$ws.Cells.Item(93, 5).Value = 5945
$ws.Cells.Item(93, 7).Value = 3
$ws.Cells.Item(93, 8).Value = 170

# --- Row 98: Senegal ---
$ws.Cells.Item(98, 2).Value = 15348
$ws.Cells.Item(98, 3).Value = 17
$ws.Cells.Item(98, 4).Value = 13637
$ws.Cells.Item(98, 5).Value = 1395

# --- Row 139: Malta ---
$ws.Cells.Item(139, 2).Value = 4160
$ws.Cells.Item(139, 3).Value = 112
$ws.Cells.Item(139, 4).Value = 3106
$ws.Cells.Item(139, 5).Value = 1009
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = 45

# --- Rows 176/177/178: Gibraltar jumps ahead of Taiwan and Burundi with fresh data ---
$ws.Cells.Item(176, 1).Value = "Gibraltar"
$ws.Cells.Item(176, 3).Value = 15
$ws.Cells.Item(176, 4).Value = 439
$ws.Cells.Item(176, 5).Value = 92
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = "Taiwan"
$ws.Cells.Item(177, 2).Value = 531
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(177, 4).Value = 491
$ws.Cells.Item(177, 5).Value = 33
$ws.Cells.Item(177, 8).Value = 7

$ws.Cells.Item(178, 1).Value = "Burundi"
$ws.Cells.Item(178, 2).Value = 529
$ws.Cells.Item(178, 4).Value = 497
$ws.Cells.Item(178, 5).Value = 31
$ws.Cells.Item(178, 8).Value = 1
